$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.857.77"
$ws.Range("E2").Value = "  +2.76%  "
$ws.Range("D3").Value = "2.609.00"
$ws.Range("E3").Value = "  +1.51%  "
$ws.Range("D5").Value = "'580.47"
$ws.Range("E5").Value = "  +4.98%  "
$ws.Range("D6").Value = "'144.10"
$ws.Range("E6").Value = "  +2.01%  "
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").Value = "2.634.07"
$ws.Range("E9").Value = "  +2.38%  "
$ws.Range("D10").Value = "'6.52"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").Value = "'0.106"
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("D13").Value = "'0.372"
$ws.Range("E13").Value = "  +6.15%  "
$ws.Range("D14").Value = "3.081.11"
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").Value = "60.850.39"
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "'23.59"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "'0.0000142"
$ws.Range("E17").Value = "  +4.08%  "
$ws.Range("D18").Value = "2.626.69"
$ws.Range("E18").Value = "  +2.22%  "
$ws.Range("D19").Value = "'11.33"
$ws.Range("E19").Value = "  +10.27%  "
$ws.Range("D20").Value = "'4.69"
$ws.Range("E20").Value = "  +3.66%  "
$ws.Range("D21").Value = "'350.24"
$ws.Range("E21").Value = "  +3.86%  "
$ws.Range("D22").Value = "'6.95"
$ws.Range("E22").Value = "  +8.00%  "
$ws.Range("D23").Value = "'0.997"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  +8.75%  "
$ws.Range("D25").Value = "'63.28"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  +0.87%  "
$ws.Range("D28").Value = "'7.90"
$ws.Range("E28").Value = "  +7.33%  "
$ws.Range("D29").Value = "0.0₃0801"
$ws.Range("E29").Value = "  +4.40%  "
$ws.Range("E30").Value = "  +9.76%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.37"
$ws.Range("E31").Value = "  +3.78%  "
$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "'0.997"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("D33").Value = "'163.58"
$ws.Range("E33").Value = "  +3.22%  "
$ws.Range("D34").Value = "'19.54"
$ws.Range("E34").Value = "  +2.92%  "
$ws.Range("E35").Value = "  +14.86%  "
$ws.Range("D36").Value = "'4.28"
$ws.Range("E36").Value = "  +5.47%  "
$ws.Range("E37").Value = "  +6.38%  "
$ws.Range("E38").Value = "  +10.85%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'37.97"
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'3.91"
$ws.Range("E40").Value = "  +7.03%  "
$ws.Range("D41").Value = "'308.31"
$ws.Range("E41").Value = "  +6.88%  "
$ws.Range("D42").Value = "'0.845"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "'135.22"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'0.995"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0985"
$ws.Range("E45").Value = "  +1.67%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'5.04"
$ws.Range("E46").Value = "  +11.99%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'19.86"
$ws.Range("E47").Value = "  +5.93%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.608"
$ws.Range("E48").Value = "  +2.94%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "'0.0553"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'20.33"
$ws.Range("E50").Value = "  +9.54%  "
$ws.Range("D51").Value = "'0.0243"
$ws.Range("E51").Value = "  +4.29%  "
